$wb = $excel.ActiveWorkbook

# --- lsh_unit_categories: fix row 23 (B/C columns) and append a new row 26 ---
$wsUnit = $wb.Worksheets.Item("lsh_unit_categories")

$wsUnit.Cells.Item(23, 2).Value = "Bráðamóttaka"
$wsUnit.Cells.Item(23, 3).Value = "emergency_room"

$wsUnit.Cells.Item(26, 1).Value = "Fv-G3 BM Göngu"
$wsUnit.Cells.Item(26, 2).Value = "Bráðamóttaka"
$wsUnit.Cells.Item(26, 3).Value = "emergency_room"
$wsUnit.Cells.Item(26, 4).Value = "home"
$wsUnit.Cells.Item(26, 5).Value = 1

$wsUnit.Range("A8").Select()
$wsUnit.Range("B26").Select()

# --- lsh_sheet_names: append a new row 14 ---
$wsNames = $wb.Worksheets.Item("lsh_sheet_names")
$wsNames.Cells.Item(14, 1).Value = "Takmörkun meðferðar"

$wsNames.Range("D18").Select()

# --- make lsh_sheet_names the active (tab-selected) sheet ---
$wsNames.Activate()
